$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
try {
  $v = $ws.StandardHeight
  Write-Output "StandardHeight: $v"
} catch {
  Write-Output "StandardHeight failed: $_"
}
try {
  $ws.StandardHeight = 15
  Write-Output "set ok"
} catch {
  Write-Output "set failed: $_"
}
